$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1567.7368
$ws.Range("I17").Value = 449.25
$ws.Range("J17").Value = 1866
$ws.Range("K17").Value = 1347.75
$ws.Range("L17").Value = 5598
$ws.Range("M17").Value = -1179.75
$ws.Range("N17").Value = -5934
$ws.Range("H121").Value = 991.25
$ws.Range("I121").Value = 550
$ws.Range("J121").Value = 1054.2858
$ws.Range("K121").Value = 1650
$ws.Range("L121").Value = 3162.8574
$ws.Range("M121").Value = 97
$ws.Range("N121").Value = -6656.857400000001
$ws.Range("H123").Value = 45062.582
$ws.Range("I123").Value = 9000
$ws.Range("J123").Value = 47549.656
$ws.Range("K123").Value = 9000
$ws.Range("L123").Value = 47549.656
$ws.Range("M123").Value = -4100
$ws.Range("N123").Value = -57349.656
$ws.Range("H130").Value = 41143.637
$ws.Range("J130").Value = 41143.637
$ws.Range("L130").Value = 41143.637
$ws.Range("N130").Value = -51183.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20198.463
$ws.Range("I32").Value = 5546.6626
$ws.Range("K32").Value = 5546.6626
$ws.Range("M32").Value = -5259.6626
$ws.Range("H131").Value = 46413.855
$ws.Range("J131").Value = 46413.855
$ws.Range("L131").Value = 46413.855
$ws.Range("N131").Value = -56493.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 47222
$ws.Range("J122").Value = 47222
$ws.Range("L122").Value = 47222
$ws.Range("N122").Value = -57022
$ws.Range("H126").Value = 34180
$ws.Range("J126").Value = 34180
$ws.Range("L126").Value = 34180
$ws.Range("N126").Value = -44060
$ws.Range("H129").Value = 49986.332
$ws.Range("J129").Value = 49986.332
$ws.Range("L129").Value = 49986.332
$ws.Range("N129").Value = -59986.332
$ws.Range("H130").Value = 48688.89
$ws.Range("J130").Value = 48688.89
$ws.Range("L130").Value = 48688.89
$ws.Range("N130").Value = -58728.89
$ws.Range("H135").Value = 44190.76
$ws.Range("J135").Value = 44190.76
$ws.Range("L135").Value = 44190.76
$ws.Range("N135").Value = -54330.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45256.5
$ws.Range("J20").Value = 45256.5
$ws.Range("L20").Value = 45256.5
$ws.Range("N20").Value = -45728.5
$ws.Range("H30").Value = 45256.5
$ws.Range("J30").Value = 45256.5
$ws.Range("L30").Value = 45256.5
$ws.Range("N30").Value = -45438.5
$ws.Range("H64").Value = 32250
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 32250
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 32250
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -32746
$ws.Range("H67").Value = 32250
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 32250
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 32250
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -33966
$ws.Range("H99").Value = 1624.8422
$ws.Range("I99").Value = 1499.3846
$ws.Range("J99").Value = 1896.6666
$ws.Range("K99").Value = 1499.3846
$ws.Range("L99").Value = 1896.6666
$ws.Range("M99").Value = -1.384600000000091
$ws.Range("N99").Value = -4892.6666
$ws.Range("H124").Value = 20076.5
$ws.Range("J124").Value = 20076.5
$ws.Range("L124").Value = 20076.5
$ws.Range("N124").Value = -24986.5
$ws.Range("H126").Value = 1624.8422
$ws.Range("I126").Value = 1499.3846
$ws.Range("J126").Value = 1896.6666
$ws.Range("K126").Value = 4498.1538
$ws.Range("L126").Value = 5689.9998
$ws.Range("M126").Value = -2028.1538
$ws.Range("N126").Value = -10629.9998
$ws.Range("H128").Value = 45256.5
$ws.Range("J128").Value = 45256.5
$ws.Range("L128").Value = 45256.5
$ws.Range("N128").Value = -55216.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 7324.3125
$ws.Range("I70").Value = 7799.4
$ws.Range("J70").Value = 7108.364
$ws.Range("K70").Value = 23398.2
$ws.Range("L70").Value = 21325.092
$ws.Range("M70").Value = -23083.2
$ws.Range("N70").Value = -21955.092
$ws.Range("H73").Value = 7324.3125
$ws.Range("I73").Value = 7799.4
$ws.Range("J73").Value = 7108.364
$ws.Range("K73").Value = 23398.2
$ws.Range("L73").Value = 21325.092
$ws.Range("M73").Value = -22306.2
$ws.Range("N73").Value = -23509.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1316.4517
$ws.Range("I102").Value = 971.84
$ws.Range("J102").Value = 2752.3333
$ws.Range("K102").Value = 971.84
$ws.Range("L102").Value = 2752.3333
$ws.Range("M102").Value = 650.16
$ws.Range("N102").Value = -5996.3333
$ws.Range("H132").Value = 2429.2258
$ws.Range("I132").Value = 2229.8
$ws.Range("J132").Value = 2791.818
$ws.Range("K132").Value = 6689.400000000001
$ws.Range("L132").Value = 8375.454000000002
$ws.Range("M132").Value = -4159.400000000001
$ws.Range("N132").Value = -13435.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 35150
$ws.Range("J64").Value = 35150
$ws.Range("L64").Value = 35150
$ws.Range("N64").Value = -35600
$ws.Range("H67").Value = 35150
$ws.Range("J67").Value = 35150
$ws.Range("L67").Value = 35150
$ws.Range("N67").Value = -36710
$ws.Range("H130").Value = 55564
$ws.Range("J130").Value = 55564
$ws.Range("L130").Value = 55564
$ws.Range("N130").Value = -65604
$ws.Range("H131").Value = 36045
$ws.Range("J131").Value = 36045
$ws.Range("L131").Value = 36045
$ws.Range("N131").Value = -46125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H108").Value = 27795
$ws.Range("J108").Value = 27795
$ws.Range("L108").Value = 27795
$ws.Range("N108").Value = -35475
$ws.Range("H122").Value = 8637.069
$ws.Range("I122").Value = 10697
$ws.Range("K122").Value = 32091
$ws.Range("M122").Value = -29641
$ws.Range("H127").Value = 31139.834
$ws.Range("J127").Value = 31139.834
$ws.Range("L127").Value = 31139.834
$ws.Range("N127").Value = -41059.834
$ws.Range("H128").Value = 275357.5
$ws.Range("J128").Value = 275357.5
$ws.Range("L128").Value = 275357.5
$ws.Range("N128").Value = -285317.5
$ws.Range("H135").Value = 66096.664
$ws.Range("J135").Value = 73316
$ws.Range("L135").Value = 73316
$ws.Range("N135").Value = -83456
